$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "2024-09-03 12:18:17"
$ws.Range("B6").Value = "ORM-0515835"
$ws.Range("C6").Value = "Z195631098"
$ws.Range("D6").Value = "C18H19CL2N5OS"
$ws.Range("E6").Value = "Duplicate"
$ws.Range("F6").Value = "/home/robekott/ERAT/examples/compound_test.sdf"

$ws.Range("A7").Value = "2024-09-03 12:18:18"
$ws.Range("B7").Value = "ORM-0515836"
$ws.Range("C7").Value = "Z2754556176"
$ws.Range("D7").Value = "C17H28N4O2"
$ws.Range("E7").Value = "Duplicate"
$ws.Range("F7").Value = "/home/robekott/ERAT/examples/compound_test.sdf"
